$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Find the first repeating element in an array of integers"
$ws.Range("B18").Value = "FirstRepeatingElement"

$ws.Range("A14").Select()
